$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28").ClearContents()
$ws.Range("N28").ClearContents()

$ws.Range("H40").Value = 4559.8667
$ws.Range("J40").Value = 6134.1665
$ws.Range("L40").Value = 6134.1665
$ws.Range("N40").Value = -6484.1665

$ws.Range("H44").Value = 65050
$ws.Range("J44").Value = 65050
$ws.Range("L44").Value = 65050
$ws.Range("N44").Value = -65974

$ws.Range("H62").Value = 3922.2354
$ws.Range("J62").Value = 4496.6
$ws.Range("L62").Value = 4496.6
$ws.Range("N62").Value = -5744.6

$ws.Range("H65").Value = 3922.2354
$ws.Range("J65").Value = 4496.6
$ws.Range("L65").Value = 22483
$ws.Range("N65").Value = -28723

$ws.Range("H86").Value = 4943.8125
$ws.Range("I86").Value = 5004
$ws.Range("K86").Value = 5004
$ws.Range("M86").Value = -3881

$ws.Range("H89").Value = 4943.8125
$ws.Range("I89").Value = 5004
$ws.Range("K89").Value = 25020
$ws.Range("M89").Value = -19404

$ws.Range("H132").Value = 3123.0193
$ws.Range("I132").Value = 3198.5334
$ws.Range("J132").Value = 2637.5715
$ws.Range("K132").Value = 9595.600199999999
$ws.Range("L132").Value = 7912.7145
$ws.Range("M132").Value = -7065.600199999999
$ws.Range("N132").Value = -12972.7145

$ws.Range("H133").Value = 107994.5
$ws.Range("J133").Value = 107994.5
$ws.Range("L133").Value = 107994.5
$ws.Range("N133").Value = -118114.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1695304.9
$ws.Range("I32").Value = 844081.5
$ws.Range("K32").Value = 844081.5
$ws.Range("M32").Value = -843794.5

$ws.Range("H61").Value = 7500

$ws.Range("H74").Value = 18975058
$ws.Range("I74").Value = 170175.33
$ws.Range("J74").Value = 50003110
$ws.Range("K74").Value = 170175.33
$ws.Range("L74").Value = 50003110
$ws.Range("M74").Value = -169301.33
$ws.Range("N74").Value = -50004858

$ws.Range("H77").Value = 18975058
$ws.Range("I77").Value = 170175.33
$ws.Range("J77").Value = 50003110
$ws.Range("K77").Value = 850876.6499999999
$ws.Range("L77").Value = 250015550
$ws.Range("M77").Value = -846508.6499999999
$ws.Range("N77").Value = -250024286

$ws.Range("H122").Value = 26317468
$ws.Range("I122").Value = 32259730
$ws.Range("J122").Value = 1734.2858
$ws.Range("K122").Value = 96779190
$ws.Range("L122").Value = 5202.857400000001
$ws.Range("M122").Value = -96776740
$ws.Range("N122").Value = -10102.8574

$ws.Range("H132").Value = 3848595
$ws.Range("I132").Value = 6994826
$ws.Range("K132").Value = 20984478
$ws.Range("M132").Value = -20981948

$ws.Range("H136").Value = 7500

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5047.1113
$ws.Range("J86").Value = 5121.5
$ws.Range("L86").Value = 5121.5
$ws.Range("N86").Value = -7367.5

$ws.Range("H89").Value = 5047.1113
$ws.Range("J89").Value = 5121.5
$ws.Range("L89").Value = 25607.5
$ws.Range("N89").Value = -36839.5

$ws.Range("H96").Value = 40000
$ws.Range("J96").Value = 40000
$ws.Range("L96").Value = 40000
$ws.Range("N96").Value = -45492

$ws.Range("H134").Value = 1948.25
$ws.Range("J134").Value = 1898.5
$ws.Range("L134").Value = 5695.5
$ws.Range("N134").Value = -10765.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2121374.2
$ws.Range("I31").Value = 1219.5333
$ws.Range("K31").Value = 1219.5333
$ws.Range("M31").Value = -924.5333000000001

$ws.Range("H34").Value = 2121374.2
$ws.Range("I34").Value = 1219.5333
$ws.Range("K34").Value = 1219.5333
$ws.Range("M34").Value = -1017.5333

$ws.Range("H86").Value = 11506.267
$ws.Range("J86").Value = 15705.556
$ws.Range("L86").Value = 15705.556
$ws.Range("N86").Value = -17951.556

$ws.Range("H89").Value = 11506.267
$ws.Range("J89").Value = 15705.556
$ws.Range("L89").Value = 78527.78
$ws.Range("N89").Value = -89759.78

$ws.Range("H132").Value = 4577.381
$ws.Range("I132").Value = 4527.846
$ws.Range("J132").Value = 4657.875
$ws.Range("K132").Value = 13583.538
$ws.Range("L132").Value = 13973.625
$ws.Range("M132").Value = -11053.538
$ws.Range("N132").Value = -19033.625

$ws.Range("H133").Value = 266725330
$ws.Range("J133").Value = 266725330
$ws.Range("L133").Value = 266725330
$ws.Range("N133").Value = -266730390

$ws.Range("H134").Value = 5926.364
$ws.Range("I134").Value = 6354.4443
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 19063.3329
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -16528.3329
$ws.Range("N134").Value = -17070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 505.8889
$ws.Range("I5").Value = 386.0909
$ws.Range("K5").Value = 1158.2727
$ws.Range("M5").Value = -1046.2727

$ws.Range("H55").Value = 3662.1667
$ws.Range("J55").Value = 4374.75
$ws.Range("L55").Value = 13124.25
$ws.Range("N55").Value = -13478.25

$ws.Range("H81").Value = 2559.75
$ws.Range("I81").Value = 612.5
$ws.Range("J81").Value = 4507
$ws.Range("K81").Value = 1837.5
$ws.Range("L81").Value = 13521
$ws.Range("M81").Value = -714.5
$ws.Range("N81").Value = -15767

$ws.Range("H84").Value = 2559.75
$ws.Range("I84").Value = 612.5
$ws.Range("J84").Value = 4507
$ws.Range("K84").Value = 5512.5
$ws.Range("L84").Value = 40563
$ws.Range("M84").Value = 103.5
$ws.Range("N84").Value = -51795

$ws.Range("H113").Value = 622.5
$ws.Range("J113").Value = 660.53845
$ws.Range("L113").Value = 1981.61535
$ws.Range("N113").Value = -6321.61535

$ws.Range("H132").Value = 4912.5283
$ws.Range("I132").Value = 5020.3335
$ws.Range("J132").Value = 4880.9756
$ws.Range("K132").Value = 45183.0015
$ws.Range("L132").Value = 43928.7804
$ws.Range("M132").Value = -42653.0015
$ws.Range("N132").Value = -48988.7804

$ws.Range("H135").Value = 505.8889
$ws.Range("I135").Value = 386.0909
$ws.Range("K135").Value = 3474.8181
$ws.Range("M135").Value = -939.8181

$ws.Range("H137").Value = 5477.6924
$ws.Range("I137").Value = 6581.4
$ws.Range("J137").Value = 1798.6666
$ws.Range("K137").Value = 19744.2
$ws.Range("L137").Value = 5395.9998
$ws.Range("M137").Value = -14644.2
$ws.Range("N137").Value = -15595.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 55559980
$ws.Range("I126").Value = 100003980
$ws.Range("K126").Value = 300011940
$ws.Range("M126").Value = -300009470

$ws.Range("H132").Value = 2238.5945
$ws.Range("I132").Value = 2354.182
$ws.Range("K132").Value = 7062.545999999999
$ws.Range("M132").Value = -4532.545999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 142857740
$ws.Range("I22").Value = 200000480
$ws.Range("J22").Value = 896
$ws.Range("K22").Value = 200000480
$ws.Range("L22").Value = 896
$ws.Range("M22").Value = -200000185
$ws.Range("N22").Value = -1486

$ws.Range("H27").Value = 142857740
$ws.Range("I27").Value = 200000480
$ws.Range("J27").Value = 896
$ws.Range("K27").Value = 200000480
$ws.Range("L27").Value = 896
$ws.Range("M27").Value = -200000373
$ws.Range("N27").Value = -1110

$ws.Range("H40").Value = 50108.31
$ws.Range("I40").Value = 69642.836
$ws.Range("K40").Value = 69642.836
$ws.Range("M40").Value = -69506.836

$ws.Range("H82").Value = 1364.2222
$ws.Range("I82").Value = 1341.4546
$ws.Range("J82").Value = 1400
$ws.Range("K82").Value = 1341.4546
$ws.Range("L82").Value = 1400
$ws.Range("M82").Value = -980.4546
$ws.Range("N82").Value = -2122

$ws.Range("H85").Value = 1364.2222
$ws.Range("I85").Value = 1341.4546
$ws.Range("J85").Value = 1400
$ws.Range("K85").Value = 1341.4546
$ws.Range("L85").Value = 1400
$ws.Range("M85").Value = -93.45460000000003
$ws.Range("N85").Value = -3896

$ws.Range("H136").Value = 6123.6
$ws.Range("I136").Value = 4614.1113
$ws.Range("J136").Value = 6972.6875
$ws.Range("K136").Value = 13842.3339
$ws.Range("L136").Value = 20918.0625
$ws.Range("M136").Value = -11292.3339
$ws.Range("N136").Value = -26018.0625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 585
$ws.Range("I107").Value = 366.66666
$ws.Range("J107").Value = 748.75
$ws.Range("K107").Value = 1099.99998
$ws.Range("L107").Value = 2246.25
$ws.Range("M107").Value = 820.0000199999999
$ws.Range("N107").Value = -6086.25

$ws.Range("H113").Value = 1033.9231
$ws.Range("I113").Value = 1090.9166
$ws.Range("K113").Value = 3272.7498
$ws.Range("M113").Value = -1102.7498

$ws.Range("H124").Value = 148991.2
$ws.Range("J124").Value = 148991.2
$ws.Range("L124").Value = 148991.2
$ws.Range("N124").Value = -158811.2

$ws.Range("H132").Value = 4654.857
$ws.Range("I132").Value = 4317
$ws.Range("J132").Value = 5499.5
$ws.Range("K132").Value = 12951
$ws.Range("L132").Value = 16498.5
$ws.Range("M132").Value = -10421
$ws.Range("N132").Value = -21558.5
